$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split "name" into name_es (A) / name_en (B), and "description" into
# description_es (D) / description_en (E) by inserting a fresh column
# right after each, shifting everything else over.
$ws.Columns("B:B").Insert()
$ws.Columns("E:E").Insert()

# Row 1 - headers
$ws.Range("A1").Value = "name_es"
$ws.Range("B1").Value = "name_en"
$ws.Range("D1").Value = "description_es"
$ws.Range("E1").Value = "description_en"

# Row 2 - first data row
$ws.Range("A2").Value = "nombre"
$ws.Range("B2").Value = "name"
$ws.Range("D2").Value = "Test"
$ws.Range("E2").Value = "Test"

# Row 3 - second data row
$ws.Range("A3").Value = "nombre2"
$ws.Range("B3").Value = "name2"
$ws.Range("D3").Value = "Test"
$ws.Range("E3").Value = "Test"
